# Generate Report for handoff
#
# Updates the "Latest Handoff Datetime" cell (column D, row 6 — the
# ac015e3f-5436-4457-8cb1-c2fb9993a1fe.md entry which is "Ready for
# handoff") on both the "zh-cn" and "de-de" localization-status sheets
# with the freshly generated handoff timestamps.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D6").Value = "2016-01-15 10:00:23"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D6").Value = "2016-01-15 10:00:34"
